# Generate Report for Handback
# Updates the localization-status workbook to reflect that the file
# "8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md" has been handed back (in sync
# with en-US) for both the zh-cn and de-de locales: its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns get populated for the first time.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$srcMdName  = "8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
$srcMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack   # zh-cn column
$overview.Range("C2").Value = $statusHandedBack   # de-de column

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcnXlfName = "8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"
$zhcnXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44b07ec9398facf38268d5b6e7c84afdc4ccbe6f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"

$zhcn.Range("C2").Value = $statusHandedBack

$zhcn.Range("F2").Value = $srcMdName
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $srcMdUrl, "", "", $srcMdName)
$zhcn.Range("F2").Font.Underline = $true
$zhcn.Range("F2").Font.Color = 15570276

$zhcn.Range("G2").Value = $zhcnXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, "", "", $zhcnXlfName)
$zhcn.Range("G2").Font.Underline = $true
$zhcn.Range("G2").Font.Color = 15570276

$zhcn.Range("H2").Value = "2016-03-12 02:25:54"

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dedeXlfName = "8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"
$dedeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc2c60b3104014aef3802feb71c06be339879ccf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"

$dede.Range("C2").Value = $statusHandedBack

$dede.Range("F2").Value = $srcMdName
$dede.Hyperlinks.Add($dede.Range("F2"), $srcMdUrl, "", "", $srcMdName)
$dede.Range("F2").Font.Underline = $true
$dede.Range("F2").Font.Color = 15570276

$dede.Range("G2").Value = $dedeXlfName
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, "", "", $dedeXlfName)
$dede.Range("G2").Font.Underline = $true
$dede.Range("G2").Font.Color = 15570276

$dede.Range("H2").Value = "2016-03-12 02:26:00"
